# Update the cryptocurrency price/volume snapshot to the latest scrape.
# Commit: "Updated cryptos list on Sat Mar 25 11:32:40 UTC 2023 with GitHub Actions"
#
# Source data stores the Price column (D) as literal text (values like
# "1.001" or "27.524.36" are NOT numbers - they're formatted price strings),
# so we force Text number-format on that column while writing the new
# values, then clear the format tweak again so the cells end up back at the
# workbook's default (unstyled) cell format - only their text content
# changes, matching the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.521.70"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.750.06"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "324.53"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4471"
$ws.Range("E7").Value = "  +4.04%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "0.07495"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "41.98"
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "20.64"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "6.026"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  -3.12%  "
$ws.Range("D16").Value = "1.750.48"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "93.28"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "0.06382"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "5.854"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "27.567.82"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "2.081"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("D26").Value = "162.07"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "20.52"
$ws.Range("D28").Value = "1.950.41"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").Value = "2.087"
$ws.Range("E29").Value = "  -4.51%  "
$ws.Range("D30").Value = "125.65"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "1.080"
$ws.Range("E31").Value = "  -7.73%  "
$ws.Range("D32").Value = "3.661"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").Value = "0.09038"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "5.549"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -5.76%  "
$ws.Range("D36").Value = "0.02297"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").Value = "0.06013"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.6354"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2082"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "4.940"
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("D41").Value = "1.204"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").Value = "1.383"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").Value = "7.768"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "13.24"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "3.725"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "0.5888"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").Value = "121.98"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "1.148"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").Value = "0.06856"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").Value = "72.07"
$ws.Range("E51").Value = "  -3.80%  "

$priceRange.ClearFormats()
